$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to remain Text so numeric-looking strings
# (e.g. "1.001", "216.72") are not auto-coerced into numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$data = @(
    ,@(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '26.429.72', '  +1.32%  ')
    ,@(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.677.40', '  +2.48%  ')
    ,@(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.001', '  -0.06%  ')
    ,@(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '216.72', '  +1.29%  ')
    ,@(6, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.5310', '  +1.19%  ')
    ,@(7, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.002', '  -0.07%  ')
    ,@(8, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.2694', '  +3.56%  ')
    ,@(9, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.06395', '  +1.46%  ')
    ,@(10, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '21.70', '  +4.89%  ')
    ,@(11, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.07816', '  +1.98%  ')
    ,@(12, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.689.41', '  +2.52%  ')
    ,@(13, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '4.509', '  +1.98%  ')
    ,@(14, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.5559', '  +0.95%  ')
    ,@(15, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.0₅8323', '  +2.31%  ')
    ,@(16, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '65.60', '  +0.70%  ')
    ,@(17, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '26.475.58', '  +1.54%  ')
    ,@(18, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.001', '  -0.12%  ')
    ,@(19, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '4.736', '  +0.88%  ')
    ,@(20, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '193.44', '  +2.61%  ')
    ,@(21, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '10.30', '  +1.54%  ')
    ,@(22, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '6.343', '  +3.00%  ')
    ,@(23, 'BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '1.002', '  -0.05%  ')
    ,@(24, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '142.40', '  -2.67%  ')
    ,@(25, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.1284', '  +5.47%  ')
    ,@(26, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '7.400', '  -0.18%  ')
    ,@(27, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '16.22', '  +2.37%  ')
    ,@(28, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '1.429', '  +1.47%  ')
    ,@(29, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.06244', '  +4.79%  ')
    ,@(30, 'PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '1.274', '  +1.22%  ')
    ,@(31, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '3.604', '  +4.64%  ')
    ,@(32, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '3.445', '  +0.86%  ')
    ,@(33, 'LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '1.677', '  +2.17%  ')
    ,@(34, 'ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '1.007', '  +1.94%  ')
    ,@(35, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '0.6123', '  +6.88%  ')
    ,@(36, 'HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '2.428', '  +1.33%  ')
    ,@(37, 'MXToken', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx', '2.783', '  +0.73%  ')
    ,@(38, 'FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '6.161', '  +8.35%  ')
    ,@(39, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.01629', '  +0.82%  ')
    ,@(40, 'Maker', 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr', '1.083.12', '  +4.03%  ')
    ,@(41, 'TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '0.8634', '  +1.01%  ')
    ,@(42, 'PaxDollar', 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp', '1.000', '  -0.14%  ')
    ,@(43, 'Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '100.17', '  -0.53%  ')
    ,@(44, 'RocketPoolETH', 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth', '1.822.28', '  +1.99%  ')
    ,@(45, 'Aave', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave', '57.10', '  +3.15%  ')
    ,@(46, 'EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '8.165', '  +1.48%  ')
    ,@(47, 'Frax', 'https://coinranking.com/coin/KfWtaeV1W+frax-frax', '0.9962', '  -0.31%  ')
    ,@(48, 'BabyDogeCoin', 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge', '0.0₈103', '  -3.76%  ')
    ,@(49, 'Cronos', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro', '0.05206', '  +0.62%  ')
    ,@(50, 'RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '1.472', '  +6.02%  ')
    ,@(51, 'Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '6.020', '  +1.85%  ')
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
}

Write-Output "Done updating $($data.Count) rows"
